$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell updates coming from the refreshed crypto price/volume feed ---
# Row 48-51 also shifted: BabyDogeCoin dropped off the board and USDD joined
# at the bottom, so each of Coin/Link/Price/Volume shifts up one data row.

# 1) Updates whose new text is unambiguous as text (contains letters, "%",
#    multiple "." separators, a URL, etc.) - plain assignment keeps them as text.
$textUpdates = @{
    "D2" = '26.271.63'
    "E2" = '  -0.10%  '
    "D3" = '1.594.36'
    "E3" = '  +0.31%  '
    "E4" = '  -0.04%  '
    "E5" = '  +0.52%  '
    "E6" = '  -0.54%  '
    "E7" = '  -0.02%  '
    "E8" = '  -0.31%  '
    "E9" = '  -0.33%  '
    "E10" = '  -1.88%  '
    "E11" = '  +0.50%  '
    "D12" = '1.819.24'
    "E12" = '  +0.34%  '
    "D13" = '1.593.64'
    "E13" = '  -0.25%  '
    "E14" = '  -1.02%  '
    "E15" = '  -2.06%  '
    "D17" = '26.259.63'
    "E18" = '  -1.33%  '
    "E19" = '  +1.43%  '
    "E20" = '  -1.96%  '
    "E21" = '  -0.05%  '
    "E22" = '  +0.35%  '
    "E23" = '  +0.48%  '
    "E24" = '  -2.48%  '
    "E25" = '  +0.04%  '
    "E26" = '  -0.01%  '
    "E27" = '  -1.05%  '
    "E28" = '  +0.95%  '
    "E29" = '  -0.57%  '
    "E30" = '  -1.75%  '
    "E31" = '  +0.39%  '
    "E32" = '  -0.40%  '
    "D33" = '1.418.54'
    "E33" = '  +5.77%  '
    "E34" = '  -0.19%  '
    "E35" = '  -0.89%  '
    "E36" = '  -1.55%  '
    "E37" = '  -3.98%  '
    "E38" = '  -0.72%  '
    "E40" = '  +0.33%  '
    "E41" = '  +0.00%  '
    "E42" = '  +0.77%  '
    "E43" = '  -11.07%  '
    "E44" = '  -0.22%  '
    "D45" = '1.731.44'
    "E45" = '  +0.34%  '
    "E46" = '  -1.31%  '
    "E47" = '  -0.73%  '
    "B48" = 'RenderToken'
    "C48" = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    "E48" = '  -0.85%  '
    "B49" = 'Cronos'
    "C49" = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    "E49" = '  -0.33%  '
    "B50" = 'Algorand'
    "C50" = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    "E50" = '  -2.84%  '
    "B51" = 'USDD'
    "C51" = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
    "E51" = '  +0.11%  '
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# 2) Updates whose new text looks like a plain number (e.g. "18.98") which
#    Excel would otherwise auto-convert to a numeric cell on assignment.
#    The source cells are plain text (inlineStr) in the original workbook, so
#    force text via NumberFormat "@" before assigning, then restore the
#    "Normal" style afterwards so no stray cell formatting is left behind.
$numericLookingUpdates = @{
    "D6" = '0.499'
    "D10" = '18.98'
    "D11" = '0.0851'
    "D16" = '63.85'
    "D20" = '7.35'
    "D22" = '4.30'
    "D23" = '9.04'
    "D24" = '2.09'
    "D25" = '145.05'
    "D29" = '15.11'
    "D30" = '0.0491'
    "D32" = '3.20'
    "D37" = '0.576'
    "D39" = '0.824'
    "D40" = '5.80'
    "D43" = '0.928'
    "D46" = '60.95'
    "D47" = '87.34'
    "D48" = '1.48'
    "D49" = '0.0501'
    "D50" = '0.0953'
    "D51" = '1.00'
}
foreach ($addr in $numericLookingUpdates.Keys) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $numericLookingUpdates[$addr]
    $c.Style = "Normal"
}
